$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 345
$ws.Range("I4").Value = 340
$ws.Range("K4").Value = 340
$ws.Range("M4").Value = -226

$ws.Range("H17").Value = 1875.5
$ws.Range("J17").Value = 1875.5
$ws.Range("L17").Value = 5626.5
$ws.Range("N17").Value = -5962.5

$ws.Range("H19").Value = 1500
$ws.Range("I19").Value = 1500
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1325
$ws.Range("N19").ClearContents()

$ws.Range("H28").Value = 1336.76
$ws.Range("I28").Value = 1300.2273
$ws.Range("K28").Value = 1300.2273
$ws.Range("M28").Value = -815.2273

$ws.Range("H32").Value = 2538.5881
$ws.Range("I32").Value = 1707.8889
$ws.Range("K32").Value = 1707.8889
$ws.Range("M32").Value = -1381.8889

$ws.Range("H98").Value = 2088.0833
$ws.Range("I98").Value = 851.8889
$ws.Range("K98").Value = 851.8889
$ws.Range("M98").Value = 646.1111

$ws.Range("H107").Value = 1581.875
$ws.Range("I107").Value = 1510.1666
$ws.Range("J107").Value = 1797
$ws.Range("K107").Value = 1510.1666
$ws.Range("L107").Value = 1797
$ws.Range("M107").Value = 409.8334
$ws.Range("N107").Value = -5637

$ws.Range("H111").Value = 6885.4
$ws.Range("I111").Value = 6107
$ws.Range("K111").Value = 18321
$ws.Range("M111").Value = -15254

$ws.Range("H115").Value = 689
$ws.Range("I115").Value = 586.3333
$ws.Range("K115").Value = 1758.9999
$ws.Range("M115").Value = -191.9999

$ws.Range("H116").Value = 4042.375
$ws.Range("I116").Value = 3285.375
$ws.Range("K116").Value = 3285.375
$ws.Range("M116").Value = 156.625

$ws.Range("H122").Value = 2088.0833
$ws.Range("I122").Value = 851.8889
$ws.Range("K122").Value = 2555.6667
$ws.Range("M122").Value = -105.6667000000002

$ws.Range("H127").Value = 2102
$ws.Range("I127").Value = 2102
$ws.Range("K127").Value = 6306
$ws.Range("M127").Value = -1346

$ws.Range("H138").Value = 1309.3
$ws.Range("I138").Value = 1309.3
$ws.Range("K138").Value = 3927.9
$ws.Range("M138").Value = 1212.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 837.1429000000001
$ws.Range("I2").Value = 643.3333
$ws.Range("K2").Value = 643.3333
$ws.Range("M2").Value = -530.3333

$ws.Range("H32").Value = 1312.697
$ws.Range("I32").Value = 1197.4688
$ws.Range("K32").Value = 1197.4688
$ws.Range("M32").Value = -910.4688000000001

$ws.Range("H45").Value = 1679.7273
$ws.Range("I45").Value = 1679.7273
$ws.Range("K45").Value = 1679.7273
$ws.Range("M45").Value = -1302.7273

$ws.Range("H94").Value = 625000
$ws.Range("J94").Value = 625000
$ws.Range("L94").Value = 625000
$ws.Range("N94").Value = -626802

$ws.Range("H116").Value = 837.1429000000001
$ws.Range("I116").Value = 643.3333
$ws.Range("K116").Value = 643.3333
$ws.Range("M116").Value = 1650.6667

$ws.Range("H131").Value = 84900
$ws.Range("J131").Value = 84900
$ws.Range("L131").Value = 84900
$ws.Range("N131").Value = -94980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 837.1429000000001
$ws.Range("I3").Value = 643.3333
$ws.Range("K3").Value = 643.3333
$ws.Range("M3").Value = -529.3333

$ws.Range("H86").Value = 8274.883
$ws.Range("I86").Value = 2473
$ws.Range("K86").Value = 2473
$ws.Range("M86").Value = -1350

$ws.Range("H89").Value = 8274.883
$ws.Range("I89").Value = 2473
$ws.Range("K89").Value = 12365
$ws.Range("M89").Value = -6749

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 290
$ws.Range("I7").Value = 144.5
$ws.Range("J7").Value = 726.5
$ws.Range("K7").Value = 144.5
$ws.Range("L7").Value = 726.5
$ws.Range("M7").Value = -31.5
$ws.Range("N7").Value = -952.5

$ws.Range("H16").Value = 805.875
$ws.Range("I16").Value = 774.5
$ws.Range("K16").Value = 774.5
$ws.Range("M16").Value = -487.5

$ws.Range("H28").Value = 15381
$ws.Range("J28").Value = 15381
$ws.Range("L28").Value = 15381
$ws.Range("N28").Value = -15871

$ws.Range("H107").Value = 719.5
$ws.Range("I107").Value = 697.9231
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 697.9231
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1222.0769
$ws.Range("N107").Value = -4840

$ws.Range("H113").Value = 805.875
$ws.Range("I113").Value = 774.5
$ws.Range("K113").Value = 774.5
$ws.Range("M113").Value = 1395.5

$ws.Range("H131").Value = 36500
$ws.Range("J131").Value = 36500
$ws.Range("L131").Value = 36500
$ws.Range("N131").Value = -46580

$ws.Range("H133").Value = 49998
$ws.Range("J133").Value = 49998
$ws.Range("L133").Value = 49998
$ws.Range("N133").Value = -55058

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 139.83333
$ws.Range("I6").Value = 117.8
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 353.4
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -240.4
$ws.Range("N6").Value = -976

$ws.Range("H39").Value = 6100
$ws.Range("J39").Value = 6583.3335
$ws.Range("L39").Value = 19750.0005
$ws.Range("N39").Value = -20338.0005

$ws.Range("H56").Value = 11168
$ws.Range("I56").Value = 11168
$ws.Range("K56").Value = 11168
$ws.Range("M56").Value = -10638

$ws.Range("H86").Value = 2223.889
$ws.Range("I86").Value = 1014.3333
$ws.Range("K86").Value = 3042.9999
$ws.Range("M86").Value = -1856.9999

$ws.Range("H89").Value = 2223.889
$ws.Range("I89").Value = 1014.3333
$ws.Range("K89").Value = 9128.9997
$ws.Range("M89").Value = -3200.9997

$ws.Range("H99").Value = 1507.3334
$ws.Range("I99").Value = 1507.3334
$ws.Range("K99").Value = 4522.0002
$ws.Range("M99").Value = -2276.0002

$ws.Range("H132").Value = 1610.1111
$ws.Range("J132").Value = 1660
$ws.Range("L132").Value = 14940
$ws.Range("N132").Value = -20000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2683.7693
$ws.Range("I80").Value = 2741.8572
$ws.Range("K80").Value = 2741.8572
$ws.Range("M80").Value = -1743.8572

$ws.Range("H83").Value = 2683.7693
$ws.Range("I83").Value = 2741.8572
$ws.Range("K83").Value = 13709.286
$ws.Range("M83").Value = -8717.286

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 2974.6667
$ws.Range("J132").Value = 2799.5
$ws.Range("L132").Value = 8398.5
$ws.Range("N132").Value = -13458.5

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3603
$ws.Range("I7").Value = 3243.35
$ws.Range("K7").Value = 3243.35
$ws.Range("M7").Value = -3131.35

$ws.Range("H94").Value = 65891.75
$ws.Range("J94").Value = 65891.75
$ws.Range("L94").Value = 65891.75
$ws.Range("N94").Value = -67243.75

$ws.Range("H126").Value = 3603
$ws.Range("I126").Value = 3243.35
$ws.Range("K126").Value = 9730.049999999999
$ws.Range("M126").Value = -7260.049999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 880.7
$ws.Range("I107").Value = 871.3333
$ws.Range("K107").Value = 2613.9999
$ws.Range("M107").Value = -693.9998999999998

$ws.Range("H122").Value = 4118.25
$ws.Range("I122").Value = 4118.25
$ws.Range("K122").Value = 12354.75
$ws.Range("M122").Value = -9904.75

$ws.Range("H130").Value = 32442.5
$ws.Range("J130").Value = 32442.5
$ws.Range("L130").Value = 32442.5
$ws.Range("N130").Value = -42482.5

$ws.Range("H132").Value = 3106
$ws.Range("I132").Value = 3131.25
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 9393.75
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -6863.75
$ws.Range("N132").Value = -14075
